# Updates cryptos list prices/volume figures per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds text-formatted numbers (e.g. "65.743.71", "1.00") that
# must stay text, not get auto-coerced into numeric values by Excel. Force the
# cell format to Text before assigning, then restore the default 'Normal' style
# so no stray formatting is left behind.
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D30', 'D31', 'D32', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '65.743.71'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '3.453.99'
$ws.Range('E3').Value = '  -3.86%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '597.59'
$ws.Range('E5').Value = '  -1.45%  '
$ws.Range('D6').Value = '137.46'
$ws.Range('E6').Value = '  -7.67%  '
$ws.Range('D7').Value = '3.452.12'
$ws.Range('E7').Value = '  -3.91%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '0.494'
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('D10').Value = '7.54'
$ws.Range('E10').Value = '  -5.35%  '
$ws.Range('D11').Value = '0.123'
$ws.Range('E11').Value = '  -9.64%  '
$ws.Range('D12').Value = '0.381'
$ws.Range('E12').Value = '  -7.82%  '
$ws.Range('D13').Value = '4.035.69'
$ws.Range('E13').Value = '  -3.94%  '
$ws.Range('D14').Value = '0.0000184'
$ws.Range('E14').Value = '  -10.70%  '
$ws.Range('D15').Value = '26.63'
$ws.Range('E15').Value = '  -10.34%  '
$ws.Range('D16').Value = '3.478.52'
$ws.Range('E16').Value = '  -3.01%  '
$ws.Range('D17').Value = '65.630.68'
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('D18').Value = '0.114'
$ws.Range('E18').Value = '  -2.53%  '
$ws.Range('D19').Value = '9.88'
$ws.Range('E19').Value = '  -10.90%  '
$ws.Range('D20').Value = '5.79'
$ws.Range('D21').Value = '13.79'
$ws.Range('E21').Value = '  -7.59%  '
$ws.Range('D22').Value = '395.90'
$ws.Range('E22').Value = '  -6.66%  '
$ws.Range('D23').Value = '0.549'
$ws.Range('E23').Value = '  -10.44%  '
$ws.Range('D24').Value = '73.57'
$ws.Range('E24').Value = '  -6.01%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = '3.597.73'
$ws.Range('E26').Value = '  -3.50%  '
$ws.Range('D27').Value = '0.0000107'
$ws.Range('E27').Value = '  -10.63%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  -10.51%  '
$ws.Range('D30').Value = '2.28'
$ws.Range('E30').Value = '  -8.79%  '
$ws.Range('D31').Value = '8.23'
$ws.Range('E31').Value = '  -12.52%  '
$ws.Range('D32').Value = '3.458.63'
$ws.Range('E32').Value = '  -3.66%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E34').Value = '  -7.63%  '
$ws.Range('D35').Value = '22.99'
$ws.Range('E35').Value = '  -8.19%  '
$ws.Range('D36').Value = '173.36'
$ws.Range('E36').Value = '  -1.15%  '
$ws.Range('D37').Value = '1.22'
$ws.Range('E37').Value = '  -14.38%  '
$ws.Range('D38').Value = '6.94'
$ws.Range('E38').Value = '  -10.45%  '
$ws.Range('D39').Value = '1.53'
$ws.Range('E39').Value = '  -7.97%  '
$ws.Range('D40').Value = '4.81'
$ws.Range('E40').Value = '  -13.38%  '
$ws.Range('D41').Value = '0.0781'
$ws.Range('E41').Value = '  -8.59%  '
$ws.Range('D42').Value = '0.821'
$ws.Range('E42').Value = '  -6.84%  '
$ws.Range('D43').Value = '43.58'
$ws.Range('E43').Value = '  -5.58%  '
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').Value = '4.45'
$ws.Range('E45').Value = '  -14.33%  '
$ws.Range('D46').Value = '1.64'
$ws.Range('E46').Value = '  -11.66%  '
$ws.Range('D47').Value = '23.36'
$ws.Range('E47').Value = '  -2.44%  '
$ws.Range('D48').Value = '1.11'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('D49').Value = '6.59'
$ws.Range('E49').Value = '  -7.71%  '
$ws.Range('D50').Value = '2.13'
$ws.Range('E50').Value = '  -15.65%  '
$ws.Range('D51').Value = '2.205.02'
$ws.Range('E51').Value = '  -8.48%  '

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
